# "showall and add product finished"
# - Update existing rows 2-4 price values
# - Add 7 new product rows (5-11), copying the row-4 formatting down for column A
# - Dimension will auto-extend to A1:D11

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows ---
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 2

$ws.Range("C3").Value = 5.5
$ws.Range("D3").Value = 6.7

$ws.Range("C4").Value = 0
$ws.Range("D4").Value = 0

# --- New rows data: Id, Name, CPrice, PhPrice ---
$longF = "".PadRight(99, 'f')

$newRows = @(
    @(3, "a",    0,      0),
    @(4, "b",    5.7,    0),
    @(5, "c",    0,      0),
    @(6, "d",    10,     15),
    @(7, "e",    10.025, 15),
    @(8, $longF, 1,      2),
    @(9, "zzzz", 1,      2)
)

$r = 5
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}

# Copy column-A formatting (bold, centered, bordered) from row 4 down to the new rows
$ws.Range("A4").Copy()
$ws.Range("A5:A11").PasteSpecial(-4122)
$excel.CutCopyMode = 0
